$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the "Linked List" section header to "Linked List (Basics)"
$ws.Range("B57").Value = "Linked List (Basics)"

# 2. Insert 5 new rows at 60..64 (pushes the old "Binary Tree (basic)" block from
#    61/63 down to 66/68, matching the target layout)
$ws.Rows("60:64").Insert()

# The insert copies the neighbouring rows' column-B formatting into the two
# blank spacer rows (61 and 63); clear those stray cells so the spacer rows
# only carry the E-column format, matching the source layout.
$ws.Range("B61").Clear()
$ws.Range("B63").Clear()

# 3. Row 60: new item "Merge Two Sorted Lists" (Q21)
$ws.Range("A60").Value = 21
$ws.Range("B60").Value = "Merge Two Sorted Lists"
$ws.Range("B60").Style = $ws.Range("B59").Style
$ws.Range("C60").Value = "Easy"
$ws.Range("D60").Value = "Linked List "
$ws.Range("E60").Value = 45747
$ws.Range("E60").Style = $ws.Range("E59").Style

# 4. Row 62: new section header "Linked List (Advanced)"
$ws.Range("B62").Value = "Linked List (Advanced)"
$ws.Range("B62").Style = $ws.Range("B57").Style

# 5. Row 64: new item "Add Two Numbers" (Q2)
$ws.Range("A64").Value = 2
$ws.Range("B64").Value = "Add Two Numbers"
$ws.Range("B64").Style = $ws.Range("B59").Style
$ws.Range("C64").Value = "Medium"
$ws.Range("D64").Value = "Linked List "
$ws.Range("E64").Value = 45747
$ws.Range("E64").Style = $ws.Range("E59").Style

# 6. Selection bookkeeping to match the saved view state
$ws.Range("B60").Select()
